# Updated symbol list with latest cryptocurrency price/volume data
# Preserves each cell's original Text storage type (avoids Excel auto-converting
# numeric-looking strings into Number/Percentage cells) by entering the value with
# a leading quote-prefix (exactly like a user typing '302.01 into the cell) and then
# resetting the style back to Normal so no extra formatting/style index lingers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "302.01"
Set-TextValue $ws.Range("E2") "-0.32%"
Set-TextValue $ws.Range("D3") "37.72"
Set-TextValue $ws.Range("E3") "8.32%"
Set-TextValue $ws.Range("D4") "5.005"
Set-TextValue $ws.Range("E4") "-2.64%"
Set-TextValue $ws.Range("D5") "0.07860"
Set-TextValue $ws.Range("E5") "1.37%"
Set-TextValue $ws.Range("D6") "2.219"
Set-TextValue $ws.Range("E6") "-5.25%"
Set-TextValue $ws.Range("D7") "8.024"
Set-TextValue $ws.Range("E7") "0.12%"
Set-TextValue $ws.Range("D8") "4.009"
Set-TextValue $ws.Range("D9") "0.9086"
Set-TextValue $ws.Range("E9") "-2.13%"
Set-TextValue $ws.Range("D10") "0.1861"
Set-TextValue $ws.Range("E10") "3.52%"
Set-TextValue $ws.Range("D11") "0.09206"
Set-TextValue $ws.Range("E11") "-9.08%"
Set-TextValue $ws.Range("D12") "0.08464"
Set-TextValue $ws.Range("E12") "0.09%"
Set-TextValue $ws.Range("D13") "0.03520"
Set-TextValue $ws.Range("E13") "6.13%"
Set-TextValue $ws.Range("D14") "0.09924"
Set-TextValue $ws.Range("E14") "0.31%"
Set-TextValue $ws.Range("D15") "0.001484"
Set-TextValue $ws.Range("E15") "-0.82%"
Set-TextValue $ws.Range("D16") "0.005656"
Set-TextValue $ws.Range("E16") "-1.68%"
Set-TextValue $ws.Range("D17") "3.472"
Set-TextValue $ws.Range("E17") "0.11%"
Set-TextValue $ws.Range("E18") "-4.91%"
Set-TextValue $ws.Range("E19") "2.86%"
Set-TextValue $ws.Range("D20") "0.1308"
Set-TextValue $ws.Range("E20") "-0.21%"
Set-TextValue $ws.Range("D21") "4.804"
Set-TextValue $ws.Range("E21") "10.87%"
Set-TextValue $ws.Range("D22") "0.2203"
Set-TextValue $ws.Range("E22") "-7.73%"
Set-TextValue $ws.Range("D23") "0.04642"
Set-TextValue $ws.Range("E23") "1.68%"
Set-TextValue $ws.Range("E24") "0.85%"
Set-TextValue $ws.Range("E25") "-0.32%"
Set-TextValue $ws.Range("D26") "0.0001298"
Set-TextValue $ws.Range("E26") "-0.15%"
Set-TextValue $ws.Range("D27") "0.0004745"
Set-TextValue $ws.Range("E27") "39.87%"
Set-TextValue $ws.Range("E39") "-1.36%"
Set-TextValue $ws.Range("D40") "0.04724"
Set-TextValue $ws.Range("E40") "-0.76%"
Set-TextValue $ws.Range("D41") "0.007917"
Set-TextValue $ws.Range("E41") "2.01%"
Set-TextValue $ws.Range("D42") "0.1392"
Set-TextValue $ws.Range("E42") "-1.31%"
Set-TextValue $ws.Range("D43") "0.007658"
Set-TextValue $ws.Range("E43") "9.10%"
Set-TextValue $ws.Range("D44") "0.002197"
Set-TextValue $ws.Range("E44") "2.17%"
Set-TextValue $ws.Range("D45") "0.01018"
Set-TextValue $ws.Range("E45") "10.86%"
Set-TextValue $ws.Range("D46") "0.00005990"
Set-TextValue $ws.Range("E46") "-2.14%"
Set-TextValue $ws.Range("E47") "-0.11%"
Set-TextValue $ws.Range("E48") "218.15%"
Set-TextValue $ws.Range("D49") "0.002687"
Set-TextValue $ws.Range("E49") "34.33%"
Set-TextValue $ws.Range("E50") "-0.11%"
Set-TextValue $ws.Range("D51") "0.0001998"
Set-TextValue $ws.Range("E51") "-0.11%"
